$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update container_diameter (B4) from 10 to 15
$ws.Range("B4").Value = 15

# Force recalculation so the dependent formula in B6 (fluid_mass) updates
$excel.CalculateFullRebuild()
